$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.023.71'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.054.39'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.44'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.94'
$ws.Range('E7').Value = '  -2.02%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0774'
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.110'
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.45'
$ws.Range('E12').Value = '  -4.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.885'
$ws.Range('E13').Value = '  +5.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.354.57'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.72'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.058.09'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.16'
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.000.98'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.02'
$ws.Range('E19').Value = '  -2.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0892'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.44'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.14'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.01'
$ws.Range('E25').Value = '  +4.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.48'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -3.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.09'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.42'
$ws.Range('E29').Value = '  +12.27%  '
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.69'
$ws.Range('E32').Value = '  +2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.32'
$ws.Range('E35').Value = '  +3.50%  '
$ws.Range('E36').Value = '  +5.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0840'
$ws.Range('E37').Value = '  -6.74%  '
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.24'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.07'
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0224'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0962'
$ws.Range('E43').Value = '  -11.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.71'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.01'
$ws.Range('E45').Value = '  -4.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.304.21'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  -5.66%  '
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.245.72'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.48'
$ws.Range('E51').Value = '  +1.51%  '
